$d = $word.ActiveDocument

# ===========================================================================
# 1) Title row: "VIDEO <NAME> Version 1" -> "VIDEO P6-X-X_XXXX Version 1"
#    with colour runs (X / X / XXXX / 1 in red) and a "_GoBack" bookmark
#    wrapping the whole title text.
# ===========================================================================

$titleRng = $d.Range(0, 22)
$d.Bookmarks.Add("_GoBack", $titleRng)

# "<NAME>" -> "X" (keeps the existing red formatting of that run), then
# prepend "P6-" (non-red, merges with the preceding "VIDEO " run).
$rng = $d.Content
$rng.Find.Execute("<NAME>") | Out-Null
$nameStart = $rng.Start
$rng.Text = "X"
$rng.Collapse(1)
$rng.InsertBefore("P6-")

# After the insert, the lone red "X" sits right after "VIDEO P6-".
$xPos = $nameStart + 3
$rng = $d.Range($xPos + 1, $xPos + 1)
$rng.InsertAfter("-")
$rng = $d.Range($xPos + 2, $xPos + 2)
$rng.InsertAfter("X")
$d.Range($xPos + 2, $xPos + 3).Font.Color = 255
$rng = $d.Range($xPos + 3, $xPos + 3)
$rng.InsertAfter("_")
$rng = $d.Range($xPos + 4, $xPos + 4)
$rng.InsertAfter("XXXX")
$d.Range($xPos + 4, $xPos + 8).Font.Color = 255

# " Version 1" -> " Version " (black) + "1" (red)
$rng = $d.Content
$rng.Find.Execute(" Version 1") | Out-Null
$verStart = $rng.Start
$rng.Text = " Version "
$rng.Collapse(0)
$rng.InsertAfter("1")
$d.Range($verStart + 9, $verStart + 10).Font.Color = 255

# ===========================================================================
# 2) "SCREEN CAPTURE:" filename placeholder:
#    "<File Name>" -> "P6-X-X_XXXX_captureX.trec" (with matching colour runs)
# ===========================================================================

$rng = $d.Content
$rng.Find.Execute("<File Name>") | Out-Null
$fnStart = $rng.Start
$rng.Text = "X"
$rng.Collapse(1)
$rng.InsertBefore("P6-")

$xPos2 = $fnStart + 3
$rng = $d.Range($xPos2 + 1, $xPos2 + 1)
$rng.InsertAfter("-")
$rng = $d.Range($xPos2 + 2, $xPos2 + 2)
$rng.InsertAfter("X")
$d.Range($xPos2 + 2, $xPos2 + 3).Font.Color = 255
$rng = $d.Range($xPos2 + 3, $xPos2 + 3)
$rng.InsertAfter("_")
$rng = $d.Range($xPos2 + 4, $xPos2 + 4)
$rng.InsertAfter("XXXX")
$d.Range($xPos2 + 4, $xPos2 + 8).Font.Color = 255

# position right after "...XXXX" is $xPos2 + 8
$afterXxxx = $xPos2 + 8
$rng = $d.Range($afterXxxx, $afterXxxx)
$rng.InsertAfter("_capture")
$rng = $d.Range($afterXxxx + 8, $afterXxxx + 8)
$rng.InsertAfter("X")
$d.Range($afterXxxx + 8, $afterXxxx + 9).Font.Color = 255
$rng = $d.Range($afterXxxx + 9, $afterXxxx + 9)
$rng.InsertAfter(".trec")

# ===========================================================================
# 3) Merge the stray bookmark-split "community.cypress.c" + "om" back into a
#    single run (drops the old "_GoBack" bookmark that used to live there).
# ===========================================================================

$rng = $d.Content
$rng.Find.Execute("community.cypress.com") | Out-Null
$rng.Text = "ZZZZZZZZZZZZZZZZZZZZZZ"
$rng2 = $d.Range($rng.Start, $rng.End)
$rng2.Text = "community.cypress.com"
